$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")

# Status column (B): "Ready for handoff" -> "Handed back: in sync with en-US"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("B3").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime (G): fill in real handback timestamps
$ws.Range("G2").Value = "2016-03-08 08:40:16"
$ws.Range("G3").Value = "2016-03-08 08:40:16"

# Latest Target File (E) / Latest Handback File (F) - new hyperlinked cells
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/0281424f3113abb0767b38150c794396ee05138d/e2e/71b3328f-1310-469a-b4f3-1c9f9a59fceb.md", "", "", "71b3328f-1310-469a-b4f3-1c9f9a59fceb.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e0f5bc4ae8ff92ab9becad88b3d94aa0433163d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/71b3328f-1310-469a-b4f3-1c9f9a59fceb.79291b9eff26dc1d45b365227012e2d8c7b02d19.zh-cn.xlf", "", "", "71b3328f-1310-469a-b4f3-1c9f9a59fceb.79291b9eff26dc1d45b365227012e2d8c7b02d19.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/0281424f3113abb0767b38150c794396ee05138d/e2e/a681e2b7-4689-40bd-8737-056e4b80efae.md", "", "", "a681e2b7-4689-40bd-8737-056e4b80efae.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e0f5bc4ae8ff92ab9becad88b3d94aa0433163d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a681e2b7-4689-40bd-8737-056e4b80efae.1f9bb194ef813dae4b3676c8d2a51cdbb27044c4.zh-cn.xlf", "", "", "a681e2b7-4689-40bd-8737-056e4b80efae.1f9bb194ef813dae4b3676c8d2a51cdbb27044c4.zh-cn.xlf") | Out-Null

# match existing hyperlink look (underline + custom blue) instead of the generic theme hyperlink style
$ws.Range("E2:F3").Font.Underline = $true
$ws.Range("E2:F3").Font.Color = 15570276

# --- de-de sheet ---
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("B3").Value = "Handed back: in sync with en-US"

$ws2.Range("G2").Value = "2016-03-08 08:40:25"
$ws2.Range("G3").Value = "2016-03-08 08:40:25"

$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/0281424f3113abb0767b38150c794396ee05138d/e2e/71b3328f-1310-469a-b4f3-1c9f9a59fceb.md", "", "", "71b3328f-1310-469a-b4f3-1c9f9a59fceb.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2e5ced7ebfacb8ed483700d260752e893ccfba44/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/71b3328f-1310-469a-b4f3-1c9f9a59fceb.79291b9eff26dc1d45b365227012e2d8c7b02d19.de-de.xlf", "", "", "71b3328f-1310-469a-b4f3-1c9f9a59fceb.79291b9eff26dc1d45b365227012e2d8c7b02d19.de-de.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/0281424f3113abb0767b38150c794396ee05138d/e2e/a681e2b7-4689-40bd-8737-056e4b80efae.md", "", "", "a681e2b7-4689-40bd-8737-056e4b80efae.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2e5ced7ebfacb8ed483700d260752e893ccfba44/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a681e2b7-4689-40bd-8737-056e4b80efae.1f9bb194ef813dae4b3676c8d2a51cdbb27044c4.de-de.xlf", "", "", "a681e2b7-4689-40bd-8737-056e4b80efae.1f9bb194ef813dae4b3676c8d2a51cdbb27044c4.de-de.xlf") | Out-Null

$ws2.Range("E2:F3").Font.Underline = $true
$ws2.Range("E2:F3").Font.Color = 15570276
